$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.114.16"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "1.813.49"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.73"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4626"
$ws.Range("E7").Value = "  +5.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3757"
$ws.Range("E8").Value = "  +2.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07425"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8652"
$ws.Range("E10").Value = "  +0.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.60"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").Value = "1.815.18"
$ws.Range("E12").Value = "  +0.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.657"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.394"
$ws.Range("E14").Value = "  +2.46%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.27"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07088"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008746"
$ws.Range("E18").Value = "  +1.00%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").Value = "27.115.52"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.319"
$ws.Range("E22").Value = "  +3.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("D24").Value = "2.043.23"
$ws.Range("E24").Value = "  +0.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.924"
$ws.Range("E25").Value = "  -2.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.60"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.210"
$ws.Range("E27").Value = "  +0.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.47"
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.280"
$ws.Range("E29").Value = "  +1.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.96"
$ws.Range("E30").Value = "  -0.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08923"
$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7740"
$ws.Range("E32").Value = "  +4.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.171"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.535"
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.902"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.124"
$ws.Range("E37").Value = "  +3.06%  "

$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05234"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.274"
$ws.Range("E40").Value = "  +3.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.924"
$ws.Range("E41").Value = "  +3.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.370"
$ws.Range("E42").Value = "  +18.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5288"
$ws.Range("E43").Value = "  +0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1679"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.617"
$ws.Range("E45").Value = "  +1.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5042"
$ws.Range("E46").Value = "  +1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.42"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.04"
$ws.Range("E48").Value = "  +0.72%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.673"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06320"
$ws.Range("E51").Value = "  +0.05%  "
